$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new bullet "Außerdem: Speicherung mit allen Kennwerten in
#    .sif-Datei" right before the "Auswertung der Spektren" bullet.
#    The run is split where the (moved) "_GoBack" bookmark sits.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Auswertung der Spektren`r") {
        $p.Range.InsertParagraphBefore()
        $newPara = $d.Paragraphs.Item($i)
        $newStart = $newPara.Range.Start
        $newPara.Range.Text = "Außerdem: Speicherung mit allen Kennwerten in .sif-Datei"
        $newPara.Range.ListFormat.ListLevelNumber = 3

        # Split the run at offset 17 ("Außerdem: Speiche" | "rung...") by
        # dropping the relocated "_GoBack" bookmark there (mirrors Word's
        # "last edit position" bookmark).
        $bmPos = $newStart + 17
        $bmRange = $d.Range($bmPos, $bmPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)
        break
    }
}

# ---------------------------------------------------------------------
# 2) Insert a new bullet "Darstellung der Kennwerte über die Zeit" right
#    after "Wissenschaftliche Analyse erfolgt anhand der .dat-Datei".
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Wissenschaftliche Analyse erfolgt anhand der .dat-Datei`r") {
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = "Darstellung der Kennwerte über die Zeit"
        $newPara.Range.ListFormat.ListLevelNumber = 3
        break
    }
}

# ---------------------------------------------------------------------
# 3) Collapse the 4 runs describing the "Browse Files" button into one.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Auswählen von Dateien (Schaltfläche „Browse Files“)`r") {
        $full = $d.Range($p.Range.Start, $p.Range.End)
        $full.Text = ""
        $ins = $d.Range($p.Range.Start, $p.Range.Start)
        $ins.InsertAfter("Auswählen von Dateien (Schaltfläche „Browse Files“)")
        break
    }
}

# ---------------------------------------------------------------------
# 4) Collapse the 6 runs describing the "Calculate" button into one.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Berechnen von Kennwerten (Schaltfläche „Calculate“)`r") {
        $full = $d.Range($p.Range.Start, $p.Range.End)
        $full.Text = ""
        $ins = $d.Range($p.Range.Start, $p.Range.Start)
        $ins.InsertAfter("Berechnen von Kennwerten (Schaltfläche „Calculate“)")
        break
    }
}
